$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Punto 4"

# Update header row (D1 and E1 are new headers; copy C1's header style/format)
$ws.Range("C1").Value = "Polinomial"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "Lagrange"
$ws.Range("E1").Value = "Numpy"

# Update data values for columns C, D, E
$ws.Range("C2").Value = 777.2104157209396
$ws.Range("D2").Value = 486226.493442297
$ws.Range("E2").Value = 33.84980254112725
$ws.Range("C3").Value = 5140370.536733359
$ws.Range("D3").Value = 85222608.6949904
$ws.Range("E3").Value = -99.06767477112044
$ws.Range("C4").Value = 385.5502156019211
$ws.Range("D4").Value = 61899.97217154503
$ws.Range("E4").Value = 18.19447764223423
$ws.Range("C5").Value = -667668.143111825
$ws.Range("D5").Value = -8935969.021080256
$ws.Range("E5").Value = 230.5733589715774
$ws.Range("C6").Value = 25.73574340343475
$ws.Range("D6").Value = 2858.679813146591
$ws.Range("E6").Value = 47.97988458354939
$ws.Range("C7").Value = 34.86005675792694
$ws.Range("D7").Value = 89.99121451377869
$ws.Range("E7").Value = 20.46013477993529
$ws.Range("C8").Value = 14697.83473551273
$ws.Range("D8").Value = 171548.7018392086
$ws.Range("E8").Value = 18.88264201456604
$ws.Range("C9").Value = 8496603.760908961
$ws.Range("D9").Value = -6111427.54250741
$ws.Range("E9").Value = -4443.226321290259
$ws.Range("C10").Value = 556.9600409269333
$ws.Range("D10").Value = 3567.400134801865
$ws.Range("E10").Value = 91.961355387818
$ws.Range("C11").Value = 15335179.95163293
$ws.Range("D11").Value = 252334614.566509
$ws.Range("E11").Value = -511.2120051503416
